$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.143.18'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '1.651.17'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.40'
$ws.Range('E5').Value = '  +3.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5230'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2612'
$ws.Range('E8').Value = '  -0.43%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06333'
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.86'
$ws.Range('E10').Value = '  -1.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07688'
$ws.Range('E11').Value = '  +1.99%  '
$ws.Range('D12').Value = '1.646.84'
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.428'
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('D14').Value = '1.872.57'
$ws.Range('E14').Value = '  -1.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5590'
$ws.Range('E15').Value = '  +1.53%  '
$ws.Range('D16').Value = '0.0₅8226'
$ws.Range('E16').Value = '  +3.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.36'
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('D18').Value = '26.125.82'
$ws.Range('E18').Value = '  -0.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.751'
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '188.92'
$ws.Range('E21').Value = '  +1.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.26'
$ws.Range('E22').Value = '  -0.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.228'
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.26'
$ws.Range('E25').Value = '  -2.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.461'
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1214'
$ws.Range('E27').Value = '  -2.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.91'
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.393'
$ws.Range('E29').Value = '  +3.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05930'
$ws.Range('E30').Value = '  -7.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.269'
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.444'
$ws.Range('E32').Value = '  -2.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.415'
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.659'
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9885'
$ws.Range('E35').Value = '  -1.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.391'
$ws.Range('E36').Value = '  -0.63%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.756'
$ws.Range('E37').Value = '  +0.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5673'
$ws.Range('E38').Value = '  -5.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01622'
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8606'
$ws.Range('E40').Value = '  -0.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.791'
$ws.Range('E41').Value = '  -5.80%  '
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('D43').Value = '1.029.72'
$ws.Range('E43').Value = '  -7.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.34'
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('D45').Value = '1.799.52'
$ws.Range('E45').Value = '  -1.24%  '
$ws.Range('D46').Value = '0.0₈109'
$ws.Range('E46').Value = '  -1.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.12'
$ws.Range('E47').Value = '  +1.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05193'
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('E51').Value = '  -0.41%  '

# Row 48 becomes EnergySwap, row 49 becomes Frax (swap + new figures)
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.138'
$ws.Range('E48').Value = '  +1.05%  '
$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.003'
$ws.Range('E49').Value = '  -0.09%  '
